$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1978043.2
$ws.Range("J17").Value = 2067913.4
$ws.Range("L17").Value = 6203740.199999999
$ws.Range("N17").Value = -6204076.199999999
$ws.Range("H19").Value = 570.0714
$ws.Range("I19").Value = 590.2
$ws.Range("K19").Value = 590.2
$ws.Range("M19").Value = -415.2
$ws.Range("H43").Value = 7422.7393
$ws.Range("I43").Value = 4011.3333
$ws.Range("J43").Value = 9615.786
$ws.Range("K43").Value = 4011.3333
$ws.Range("L43").Value = 9615.786
$ws.Range("M43").Value = -3942.3333
$ws.Range("N43").Value = -9753.786
$ws.Range("H48").Value = 5013.759
$ws.Range("J48").Value = 5013.759
$ws.Range("L48").Value = 15041.277
$ws.Range("N48").Value = -15625.277
$ws.Range("H52").Value = 574.1429000000001
$ws.Range("I52").Value = 229
$ws.Range("K52").Value = 687
$ws.Range("M52").Value = -527
$ws.Range("H56").Value = 5013.759
$ws.Range("J56").Value = 5013.759
$ws.Range("L56").Value = 15041.277
$ws.Range("N56").Value = -16109.277
$ws.Range("H69").Value = 9116.6875
$ws.Range("I69").Value = 7021
$ws.Range("J69").Value = 12609.5
$ws.Range("K69").Value = 21063
$ws.Range("L69").Value = 37828.5
$ws.Range("M69").Value = -20189
$ws.Range("N69").Value = -39576.5
$ws.Range("H72").Value = 9116.6875
$ws.Range("I72").Value = 7021
$ws.Range("J72").Value = 12609.5
$ws.Range("K72").Value = 63189
$ws.Range("L72").Value = 113485.5
$ws.Range("M72").Value = -58821
$ws.Range("N72").Value = -122221.5
$ws.Range("H80").Value = 53602.152
$ws.Range("J80").Value = 3394.4736
$ws.Range("L80").Value = 10183.4208
$ws.Range("N80").Value = -12179.4208
$ws.Range("H83").Value = 53602.152
$ws.Range("J83").Value = 3394.4736
$ws.Range("L83").Value = 30550.2624
$ws.Range("N83").Value = -40534.2624
$ws.Range("H92").Value = 1024.5
$ws.Range("I92").Value = 916.1111
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 916.1111
$ws.Range("L92").Value = 2000
$ws.Range("M92").Value = 331.8889
$ws.Range("N92").Value = -4496
$ws.Range("H96").Value = 1070.75
$ws.Range("I96").Value = 594.5
$ws.Range("J96").Value = 2499.5
$ws.Range("K96").Value = 1783.5
$ws.Range("L96").Value = 7498.5
$ws.Range("M96").Value = -410.5
$ws.Range("N96").Value = -10244.5
$ws.Range("H97").Value = 2572.8462
$ws.Range("J97").Value = 2370.5833
$ws.Range("L97").Value = 7111.749899999999
$ws.Range("N97").Value = -8103.749899999999
$ws.Range("H100").Value = 66508.39
$ws.Range("I100").Value = 78325.69500000001
$ws.Range("J100").Value = 51145.9
$ws.Range("K100").Value = 78325.69500000001
$ws.Range("L100").Value = 51145.9
$ws.Range("M100").Value = -77784.69500000001
$ws.Range("N100").Value = -52227.9
$ws.Range("H116").Value = 11113900
$ws.Range("I116").Value = 15280176
$ws.Range("K116").Value = 15280176
$ws.Range("M116").Value = -15276734
$ws.Range("H121").Value = 1671.6364
$ws.Range("J121").Value = 1671.6364
$ws.Range("L121").Value = 5014.9092
$ws.Range("N121").Value = -8508.9092
$ws.Range("H132").Value = 1669956
$ws.Range("I132").Value = 3476.4375
$ws.Range("J132").Value = 8335874
$ws.Range("K132").Value = 10429.3125
$ws.Range("L132").Value = 25007622
$ws.Range("M132").Value = -7899.3125
$ws.Range("N132").Value = -25012682
$ws.Range("H136").Value = 68126.836
$ws.Range("J136").Value = 68126.836
$ws.Range("L136").Value = 68126.836
$ws.Range("N136").Value = -78326.836
$ws.Range("H137").Value = 9642.111000000001
$ws.Range("I137").Value = 12712.125
$ws.Range("J137").Value = 3502.0833
$ws.Range("K137").Value = 38136.375
$ws.Range("L137").Value = 10506.2499
$ws.Range("M137").Value = -35586.375
$ws.Range("N137").Value = -15606.2499
$ws.Range("H138").Value = 169866.28
$ws.Range("J138").Value = 3913.2415
$ws.Range("L138").Value = 11739.7245
$ws.Range("N138").Value = -22019.7245
$ws.Range("H141").Value = 8737.440000000001
$ws.Range("I141").Value = 8765.272000000001
$ws.Range("J141").Value = 8533.333000000001
$ws.Range("K141").Value = 26295.816
$ws.Range("L141").Value = 25599.999
$ws.Range("M141").Value = -21115.816
$ws.Range("N141").Value = -35959.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3652.4285
$ws.Range("I2").Value = 3981.4375
$ws.Range("K2").Value = 3981.4375
$ws.Range("M2").Value = -3868.4375
$ws.Range("H8").Value = 7714.7144
$ws.Range("J8").Value = 8917.166999999999
$ws.Range("L8").Value = 8917.166999999999
$ws.Range("N8").Value = -9205.166999999999
$ws.Range("H28").Value = 57680504
$ws.Range("I28").Value = 15569.25
$ws.Range("K28").Value = 15569.25
$ws.Range("M28").Value = -15377.25
$ws.Range("H31").Value = 5099.9165
$ws.Range("I31").Value = 4654.4546
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 4654.4546
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -4360.4546
$ws.Range("N31").Value = -10588
$ws.Range("H32").Value = 14341.514
$ws.Range("I32").Value = 14163.193
$ws.Range("K32").Value = 14163.193
$ws.Range("M32").Value = -13876.193
$ws.Range("H45").Value = 259661.5
$ws.Range("I45").Value = 510573
$ws.Range("J45").Value = 8750
$ws.Range("K45").Value = 510573
$ws.Range("L45").Value = 8750
$ws.Range("M45").Value = -510196
$ws.Range("N45").Value = -9504
$ws.Range("H61").Value = 11072.154
$ws.Range("I61").Value = 13169.294
$ws.Range("K61").Value = 13169.294
$ws.Range("M61").Value = -12957.294
$ws.Range("H82").Value = 64449
$ws.Range("J82").Value = 64449
$ws.Range("L82").Value = 64449
$ws.Range("N82").Value = -65171
$ws.Range("H85").Value = 64449
$ws.Range("J85").Value = 64449
$ws.Range("L85").Value = 64449
$ws.Range("N85").Value = -66945
$ws.Range("H88").Value = 1537.1052
$ws.Range("I88").Value = 1653
$ws.Range("K88").Value = 1653
$ws.Range("M88").Value = -1247
$ws.Range("H91").Value = 1537.1052
$ws.Range("I91").Value = 1653
$ws.Range("K91").Value = 1653
$ws.Range("M91").Value = -249
$ws.Range("H97").Value = 6901150
$ws.Range("I97").Value = 5644.905
$ws.Range("K97").Value = 5644.905
$ws.Range("M97").Value = -5148.905
$ws.Range("H99").Value = 57680504
$ws.Range("I99").Value = 15569.25
$ws.Range("K99").Value = 15569.25
$ws.Range("M99").Value = -12574.25
$ws.Range("H110").Value = 2548.8333
$ws.Range("I110").Value = 2128.1538
$ws.Range("K110").Value = 2128.1538
$ws.Range("M110").Value = -83.15380000000005
$ws.Range("H116").Value = 3652.4285
$ws.Range("I116").Value = 3981.4375
$ws.Range("K116").Value = 3981.4375
$ws.Range("M116").Value = -1687.4375
$ws.Range("H122").Value = 1005337.2
$ws.Range("I122").Value = 5327.8945
$ws.Range("K122").Value = 15983.6835
$ws.Range("M122").Value = -13533.6835
$ws.Range("H136").Value = 11072.154
$ws.Range("I136").Value = 13169.294
$ws.Range("K136").Value = 39507.882
$ws.Range("M136").Value = -36957.882

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3652.4285
$ws.Range("I3").Value = 3981.4375
$ws.Range("K3").Value = 3981.4375
$ws.Range("M3").Value = -3867.4375
$ws.Range("H22").Value = 266.66666
$ws.Range("I22").Value = 266.66666
$ws.Range("K22").Value = 266.66666
$ws.Range("M22").Value = -93.66665999999998
$ws.Range("H80").Value = 296.55
$ws.Range("J80").Value = 247
$ws.Range("L80").Value = 247
$ws.Range("N80").Value = -2243
$ws.Range("H82").Value = 21926.092
$ws.Range("I82").Value = 7658.125
$ws.Range("J82").Value = 59974
$ws.Range("K82").Value = 7658.125
$ws.Range("L82").Value = 59974
$ws.Range("M82").Value = -7275.125
$ws.Range("N82").Value = -60740
$ws.Range("H83").Value = 296.55
$ws.Range("J83").Value = 247
$ws.Range("L83").Value = 1235
$ws.Range("N83").Value = -11219
$ws.Range("H85").Value = 21926.092
$ws.Range("I85").Value = 7658.125
$ws.Range("J85").Value = 59974
$ws.Range("K85").Value = 7658.125
$ws.Range("L85").Value = 59974
$ws.Range("M85").Value = -6332.125
$ws.Range("N85").Value = -62626
$ws.Range("H86").Value = 7298.2
$ws.Range("I86").Value = 9430.833000000001
$ws.Range("K86").Value = 9430.833000000001
$ws.Range("M86").Value = -8307.833000000001
$ws.Range("H89").Value = 7298.2
$ws.Range("I89").Value = 9430.833000000001
$ws.Range("K89").Value = 47154.165
$ws.Range("M89").Value = -41538.165
$ws.Range("H94").Value = 2697.7878
$ws.Range("I94").Value = 1238.9546
$ws.Range("J94").Value = 5615.4546
$ws.Range("K94").Value = 1238.9546
$ws.Range("L94").Value = 5615.4546
$ws.Range("M94").Value = -787.9546
$ws.Range("N94").Value = -6517.4546
$ws.Range("H97").Value = 4420
$ws.Range("I97").Value = 4420
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 4420
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -3429
$ws.Range("H107").Value = 2395.7778
$ws.Range("J107").Value = 319.5
$ws.Range("L107").Value = 319.5
$ws.Range("N107").Value = -4159.5
$ws.Range("H134").Value = 2318.1724
$ws.Range("I134").Value = 1237.8182
$ws.Range("J134").Value = 5713.5713
$ws.Range("K134").Value = 3713.4546
$ws.Range("L134").Value = 17140.7139
$ws.Range("M134").Value = -1178.4546
$ws.Range("N134").Value = -22210.7139

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 226.11765
$ws.Range("I7").Value = 221.88889
$ws.Range("J7").Value = 230.875
$ws.Range("K7").Value = 221.88889
$ws.Range("L7").Value = 230.875
$ws.Range("M7").Value = -108.88889
$ws.Range("N7").Value = -456.875
$ws.Range("H31").Value = 6787.6763
$ws.Range("I31").Value = 6965.207
$ws.Range("J31").Value = 5758
$ws.Range("K31").Value = 6965.207
$ws.Range("L31").Value = 5758
$ws.Range("M31").Value = -6670.207
$ws.Range("N31").Value = -6348
$ws.Range("H34").Value = 6787.6763
$ws.Range("I34").Value = 6965.207
$ws.Range("J34").Value = 5758
$ws.Range("K34").Value = 6965.207
$ws.Range("L34").Value = 5758
$ws.Range("M34").Value = -6763.207
$ws.Range("N34").Value = -6162
$ws.Range("H58").Value = 2457.6287
$ws.Range("I58").Value = 1806.4546
$ws.Range("J58").Value = 3559.6155
$ws.Range("K58").Value = 1806.4546
$ws.Range("L58").Value = 3559.6155
$ws.Range("M58").Value = -1603.4546
$ws.Range("N58").Value = -3965.6155
$ws.Range("H62").Value = 10452.8
$ws.Range("I62").Value = 10561.75
$ws.Range("K62").Value = 10561.75
$ws.Range("M62").Value = -9937.75
$ws.Range("H65").Value = 10452.8
$ws.Range("I65").Value = 10561.75
$ws.Range("K65").Value = 52808.75
$ws.Range("M65").Value = -49688.75
$ws.Range("H88").Value = 29471.285
$ws.Range("I88").Value = 33103.668
$ws.Range("J88").Value = 26747
$ws.Range("K88").Value = 33103.668
$ws.Range("L88").Value = 26747
$ws.Range("M88").Value = -32697.668
$ws.Range("N88").Value = -27559
$ws.Range("H91").Value = 29471.285
$ws.Range("I91").Value = 33103.668
$ws.Range("J91").Value = 26747
$ws.Range("K91").Value = 33103.668
$ws.Range("L91").Value = 26747
$ws.Range("M91").Value = -31699.668
$ws.Range("N91").Value = -29555
$ws.Range("H99").Value = 3420166.2
$ws.Range("I99").Value = 5282884.5
$ws.Range("K99").Value = 5282884.5
$ws.Range("M99").Value = -5281386.5
$ws.Range("H120").Value = 66990
$ws.Range("J120").Value = 66990
$ws.Range("L120").Value = 66990
$ws.Range("N120").Value = -74248
$ws.Range("H126").Value = 3420166.2
$ws.Range("I126").Value = 5282884.5
$ws.Range("K126").Value = 15848653.5
$ws.Range("M126").Value = -15846183.5
$ws.Range("H132").Value = 1631.3235
$ws.Range("I132").Value = 1457.4839
$ws.Range("J132").Value = 3427.6667
$ws.Range("K132").Value = 4372.4517
$ws.Range("L132").Value = 10283.0001
$ws.Range("M132").Value = -1842.4517
$ws.Range("N132").Value = -15343.0001
$ws.Range("H134").Value = 3756.625
$ws.Range("J134").Value = 6908.143
$ws.Range("L134").Value = 20724.429
$ws.Range("N134").Value = -25794.429
$ws.Range("H136").Value = 2457.6287
$ws.Range("I136").Value = 1806.4546
$ws.Range("J136").Value = 3559.6155
$ws.Range("K136").Value = 5419.3638
$ws.Range("L136").Value = 10678.8465
$ws.Range("M136").Value = -2869.3638
$ws.Range("N136").Value = -15778.8465
$ws.Range("H141").Value = 178758.48
$ws.Range("J141").Value = 193348.69
$ws.Range("L141").Value = 193348.69
$ws.Range("N141").Value = -203708.69

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2055.2942
$ws.Range("I25").Value = 106.666664
$ws.Range("J25").Value = 2472.8572
$ws.Range("K25").Value = 319.999992
$ws.Range("L25").Value = 7418.571599999999
$ws.Range("M25").Value = -150.999992
$ws.Range("N25").Value = -7756.571599999999
$ws.Range("H30").Value = 2055.2942
$ws.Range("I30").Value = 106.666664
$ws.Range("J30").Value = 2472.8572
$ws.Range("K30").Value = 319.999992
$ws.Range("L30").Value = 7418.571599999999
$ws.Range("M30").Value = -217.999992
$ws.Range("N30").Value = -7622.571599999999
$ws.Range("H33").Value = 432.5
$ws.Range("I33").Value = 348.25
$ws.Range("J33").Value = 474.625
$ws.Range("K33").Value = 2089.5
$ws.Range("L33").Value = 2847.75
$ws.Range("M33").Value = -1806.5
$ws.Range("N33").Value = -3413.75
$ws.Range("H38").Value = 2180.7827
$ws.Range("J38").Value = 2580.4736
$ws.Range("L38").Value = 7741.4208
$ws.Range("N38").Value = -8435.4208
$ws.Range("H58").Value = 3226.037
$ws.Range("J58").Value = 3330.7693
$ws.Range("L58").Value = 9992.3079
$ws.Range("N58").Value = -10248.3079
$ws.Range("H102").Value = 9136.362999999999
$ws.Range("J102").Value = 10055.556
$ws.Range("L102").Value = 30166.668
$ws.Range("N102").Value = -35034.66800000001
$ws.Range("H140").Value = 437121.9
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10562.211
$ws.Range("J80").Value = 2148.8
$ws.Range("L80").Value = 2148.8
$ws.Range("N80").Value = -4144.8
$ws.Range("H83").Value = 10562.211
$ws.Range("J83").Value = 2148.8
$ws.Range("L83").Value = 10744
$ws.Range("N83").Value = -20728
$ws.Range("H107").Value = 526.4
$ws.Range("I107").Value = 563.2222
$ws.Range("J107").Value = 402.125
$ws.Range("K107").Value = 563.2222
$ws.Range("L107").Value = 402.125
$ws.Range("M107").Value = 1356.7778
$ws.Range("N107").Value = -4242.125
$ws.Range("H113").Value = 6731.9565
$ws.Range("I113").Value = 8557.6875
$ws.Range("J113").Value = 2558.8572
$ws.Range("K113").Value = 8557.6875
$ws.Range("L113").Value = 2558.8572
$ws.Range("M113").Value = -6387.6875
$ws.Range("N113").Value = -6898.8572
$ws.Range("H122").Value = 9532.091
$ws.Range("I122").Value = 6470.9414
$ws.Range("J122").Value = 19940
$ws.Range("K122").Value = 19412.8242
$ws.Range("L122").Value = 59820
$ws.Range("M122").Value = -16962.8242
$ws.Range("N122").Value = -64720
$ws.Range("H132").Value = 2424.625
$ws.Range("I132").Value = 2322.037
$ws.Range("K132").Value = 6966.110999999999
$ws.Range("M132").Value = -4436.110999999999
$ws.Range("H136").Value = 17079.35
$ws.Range("J136").Value = 17079.35
$ws.Range("L136").Value = 51238.05
$ws.Range("N136").Value = -56338.05

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 622.3333
$ws.Range("I22").Value = 622.3333
$ws.Range("K22").Value = 622.3333
$ws.Range("M22").Value = -327.3333
$ws.Range("H27").Value = 622.3333
$ws.Range("I27").Value = 622.3333
$ws.Range("K27").Value = 622.3333
$ws.Range("M27").Value = -515.3333
$ws.Range("H42").Value = 66332.336
$ws.Range("I42").Value = 66332.336
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 66332.336
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = -65769.336
$ws.Range("H46").Value = 2553.05
$ws.Range("I46").Value = 783.3333
$ws.Range("J46").Value = 3311.5
$ws.Range("K46").Value = 783.3333
$ws.Range("L46").Value = 3311.5
$ws.Range("M46").Value = -595.3333
$ws.Range("N46").Value = -3687.5
$ws.Range("H49").Value = 66332.336
$ws.Range("I49").Value = 66332.336
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 66332.336
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = -66185.336
$ws.Range("H55").Value = 1237.7222
$ws.Range("I55").Value = 452.41666
$ws.Range("J55").Value = 2808.3333
$ws.Range("K55").Value = 452.41666
$ws.Range("L55").Value = 2808.3333
$ws.Range("M55").Value = -279.41666
$ws.Range("N55").Value = -3154.3333
$ws.Range("H61").Value = 29965.723
$ws.Range("I61").Value = 2316.647
$ws.Range("K61").Value = 2316.647
$ws.Range("M61").Value = -2114.647
$ws.Range("H93").Value = 4063.9
$ws.Range("I93").Value = 4314.9375
$ws.Range("J93").Value = 3059.75
$ws.Range("K93").Value = 4314.9375
$ws.Range("L93").Value = 3059.75
$ws.Range("M93").Value = -3066.9375
$ws.Range("N93").Value = -5555.75
$ws.Range("H113").Value = 29965.723
$ws.Range("I113").Value = 2316.647
$ws.Range("K113").Value = 2316.647
$ws.Range("M113").Value = -146.6469999999999
$ws.Range("H122").Value = 5327.2285
$ws.Range("I122").Value = 4598.269
$ws.Range("J122").Value = 7433.1113
$ws.Range("K122").Value = 13794.807
$ws.Range("L122").Value = 22299.3339
$ws.Range("M122").Value = -11344.807
$ws.Range("N122").Value = -27199.3339
$ws.Range("H136").Value = 4289.9165
$ws.Range("I136").Value = 3175.1052
$ws.Range("J136").Value = 5535.8823
$ws.Range("K136").Value = 9525.3156
$ws.Range("L136").Value = 16607.6469
$ws.Range("M136").Value = -6975.3156
$ws.Range("N136").Value = -21707.6469

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("H29").Value = 9450
$ws.Range("I29").Value = 10266.667
$ws.Range("K29").Value = 10266.667
$ws.Range("M29").Value = -9976.666999999999
$ws.Range("H40").Value = 49999
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 49999
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = 49999
$ws.Range("N40").Value = -50297
$ws.Range("H62").Value = 429250.78
$ws.Range("J62").Value = 14938.25
$ws.Range("L62").Value = 14938.25
$ws.Range("N62").Value = -16186.25
$ws.Range("H64").Value = 52832.75
$ws.Range("J64").Value = 55665.5
$ws.Range("L64").Value = 55665.5
$ws.Range("N64").Value = -56161.5
$ws.Range("H65").Value = 429250.78
$ws.Range("J65").Value = 14938.25
$ws.Range("L65").Value = 74691.25
$ws.Range("N65").Value = -80931.25
$ws.Range("H67").Value = 52832.75
$ws.Range("J67").Value = 55665.5
$ws.Range("L67").Value = 55665.5
$ws.Range("N67").Value = -57381.5
$ws.Range("H86").Value = 35000
$ws.Range("J86").Value = 35000
$ws.Range("L86").Value = 35000
$ws.Range("N86").Value = -37246
$ws.Range("H89").Value = 35000
$ws.Range("J89").Value = 35000
$ws.Range("L89").Value = 175000
$ws.Range("N89").Value = -186232
$ws.Range("H122").Value = 15168.915
$ws.Range("I122").Value = 2345.147
$ws.Range("K122").Value = 7035.441
$ws.Range("M122").Value = -4585.441
$ws.Range("H132").Value = 8178.8447
$ws.Range("I132").Value = 8954.681
$ws.Range("J132").Value = 4863.909
$ws.Range("K132").Value = 26864.043
$ws.Range("L132").Value = 14591.727
$ws.Range("M132").Value = -24334.043
$ws.Range("N132").Value = -19651.727
$ws.Range("H136").Value = 292033.12
$ws.Range("I136").Value = 309309.25
$ws.Range("J136").Value = 4098
$ws.Range("K136").Value = 927927.75
$ws.Range("L136").Value = 12294
$ws.Range("M136").Value = -925377.75
$ws.Range("N136").Value = -17394
$ws.Range("M26").ClearContents()
